$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 272, shifting existing rows 272-354 down to 273-355.
$ws.Rows.Item(272).Insert()

# Populate the newly inserted row 272 with its data (same shape as surrounding rows).
$ws.Cells.Item(272, 1).Value = 4
$ws.Cells.Item(272, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(272, 3).Value = "Los Lagos"
$ws.Cells.Item(272, 4).Value = 44985
$ws.Cells.Item(272, 5).Value = 10
$ws.Cells.Item(272, 6).Value = 100112044
$ws.Cells.Item(272, 7).Value = "Perejil"
$ws.Cells.Item(272, 8).Value = "Sin especificar"
$ws.Cells.Item(272, 9).Value = "Primera"
$ws.Cells.Item(272, 10).Value = 140
$ws.Cells.Item(272, 11).Value = 7000
$ws.Cells.Item(272, 12).Value = 7000
$ws.Cells.Item(272, 13).Value = 7000
$ws.Cells.Item(272, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(272, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(272, 16).Value = 3500
$ws.Cells.Item(272, 17).Value = 2
$ws.Cells.Item(272, 18).Value = "Hortaliza"
